$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 30624
$ws.Cells.Item(95, 10).Value = 30624
$ws.Cells.Item(95, 12).Value = 30624
$ws.Cells.Item(95, 14).Value = -36116
$ws.Cells.Item(121, 8).Value = 1225.375
$ws.Cells.Item(121, 10).Value = 1225.375
$ws.Cells.Item(121, 12).Value = 3676.125
$ws.Cells.Item(121, 14).Value = -7170.125
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3197.1
$ws.Cells.Item(32, 9).Value = 3304.162
$ws.Cells.Item(32, 11).Value = 3304.162
$ws.Cells.Item(32, 13).Value = -3017.162
$ws.Cells.Item(61, 8).Value = 2624.75
$ws.Cells.Item(61, 9).Value = 2499
$ws.Cells.Item(61, 10).Value = 2666.6667
$ws.Cells.Item(61, 11).Value = 2499
$ws.Cells.Item(61, 12).Value = 2666.6667
$ws.Cells.Item(61, 13).Value = -2287
$ws.Cells.Item(61, 14).Value = -3090.6667
$ws.Cells.Item(74, 8).Value = 1758.8
$ws.Cells.Item(74, 9).Value = 1620.8889
$ws.Cells.Item(74, 11).Value = 1620.8889
$ws.Cells.Item(74, 13).Value = -746.8888999999999
$ws.Cells.Item(77, 8).Value = 1758.8
$ws.Cells.Item(77, 9).Value = 1620.8889
$ws.Cells.Item(77, 11).Value = 8104.4445
$ws.Cells.Item(77, 13).Value = -3736.4445
$ws.Cells.Item(132, 8).Value = 3280.2632
$ws.Cells.Item(132, 9).Value = 2916.2856
$ws.Cells.Item(132, 10).Value = 4299.4
$ws.Cells.Item(132, 11).Value = 8748.856800000001
$ws.Cells.Item(132, 12).Value = 12898.2
$ws.Cells.Item(132, 13).Value = -6218.856800000001
$ws.Cells.Item(132, 14).Value = -17958.2
$ws.Cells.Item(136, 8).Value = 2624.75
$ws.Cells.Item(136, 9).Value = 2499
$ws.Cells.Item(136, 10).Value = 2666.6667
$ws.Cells.Item(136, 11).Value = 7497
$ws.Cells.Item(136, 12).Value = 8000.000100000001
$ws.Cells.Item(136, 13).Value = -4947
$ws.Cells.Item(136, 14).Value = -13100.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value = 7499.5
$ws.Cells.Item(100, 10).Value = 7499.5
$ws.Cells.Item(100, 12).Value = 7499.5
$ws.Cells.Item(100, 14).Value = -9663.5
$ws.Cells.Item(107, 8).Value = 1660.909
$ws.Cells.Item(107, 9).Value = 1144.3846
$ws.Cells.Item(107, 11).Value = 1144.3846
$ws.Cells.Item(107, 13).Value = 775.6153999999999
$ws.Cells.Item(134, 8).Value = 5060.385
$ws.Cells.Item(134, 9).Value = 1198.7391
$ws.Cells.Item(134, 10).Value = 34666.332
$ws.Cells.Item(134, 11).Value = 3596.2173
$ws.Cells.Item(134, 12).Value = 103998.996
$ws.Cells.Item(134, 13).Value = -1061.2173
$ws.Cells.Item(134, 14).Value = -109068.996
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1260.6923
$ws.Cells.Item(31, 9).Value = 1199.0834
$ws.Cells.Item(31, 10).Value = 2000
$ws.Cells.Item(31, 11).Value = 1199.0834
$ws.Cells.Item(31, 12).Value = 2000
$ws.Cells.Item(31, 13).Value = -904.0834
$ws.Cells.Item(31, 14).Value = -2590
$ws.Cells.Item(34, 8).Value = 1260.6923
$ws.Cells.Item(34, 9).Value = 1199.0834
$ws.Cells.Item(34, 10).Value = 2000
$ws.Cells.Item(34, 11).Value = 1199.0834
$ws.Cells.Item(34, 12).Value = 2000
$ws.Cells.Item(34, 13).Value = -997.0834
$ws.Cells.Item(34, 14).Value = -2404
$ws.Cells.Item(62, 8).Value = 5558526
$ws.Cells.Item(62, 9).Value = 2998.1765
$ws.Cells.Item(62, 11).Value = 2998.1765
$ws.Cells.Item(62, 13).Value = -2374.1765
$ws.Cells.Item(65, 8).Value = 5558526
$ws.Cells.Item(65, 9).Value = 2998.1765
$ws.Cells.Item(65, 11).Value = 14990.8825
$ws.Cells.Item(65, 13).Value = -11870.8825
$ws.Cells.Item(107, 8).Value = 625.92
$ws.Cells.Item(107, 9).Value = 396.42856
$ws.Cells.Item(107, 11).Value = 396.42856
$ws.Cells.Item(107, 13).Value = 1523.57144
$ws.Cells.Item(132, 8).Value = 8905.588
$ws.Cells.Item(132, 9).Value = 15029.375
$ws.Cells.Item(132, 10).Value = 3462.2222
$ws.Cells.Item(132, 11).Value = 45088.125
$ws.Cells.Item(132, 12).Value = 10386.6666
$ws.Cells.Item(132, 13).Value = -42558.125
$ws.Cells.Item(132, 14).Value = -15446.6666
$ws.Cells.Item(134, 8).Value = 33335908
$ws.Cells.Item(134, 9).Value = 55558264
$ws.Cells.Item(134, 11).Value = 166674792
$ws.Cells.Item(134, 13).Value = -166672257
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 4536.1816
$ws.Cells.Item(63, 9).Value = 774.5
$ws.Cells.Item(63, 10).Value = 6685.7144
$ws.Cells.Item(63, 11).Value = 2323.5
$ws.Cells.Item(63, 12).Value = 20057.1432
$ws.Cells.Item(63, 13).Value = -1574.5
$ws.Cells.Item(63, 14).Value = -21555.1432
$ws.Cells.Item(66, 8).Value = 4536.1816
$ws.Cells.Item(66, 9).Value = 774.5
$ws.Cells.Item(66, 10).Value = 6685.7144
$ws.Cells.Item(66, 11).Value = 6970.5
$ws.Cells.Item(66, 12).Value = 60171.4296
$ws.Cells.Item(66, 13).Value = -3226.5
$ws.Cells.Item(66, 14).Value = -67659.4296
$ws.Cells.Item(129, 8).Value = 14882052
$ws.Cells.Item(129, 9).Value = 33334056
$ws.Cells.Item(129, 10).Value = 4630938.5
$ws.Cells.Item(129, 11).Value = 100002168
$ws.Cells.Item(129, 12).Value = 13892815.5
$ws.Cells.Item(129, 13).Value = -99997168
$ws.Cells.Item(129, 14).Value = -13902815.5
$ws.Cells.Item(131, 8).Value = 22728714
$ws.Cells.Item(131, 9).Value = 111111500
$ws.Cells.Item(131, 10).Value = 1709.6
$ws.Cells.Item(131, 11).Value = 333334500
$ws.Cells.Item(131, 12).Value = 5128.799999999999
$ws.Cells.Item(131, 13).Value = -333329460
$ws.Cells.Item(131, 14).Value = -15208.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 20690
$ws.Cells.Item(95, 10).Value = 20690
$ws.Cells.Item(95, 12).Value = 20690
$ws.Cells.Item(95, 14).Value = -26182
$ws.Cells.Item(122, 8).Value = 1402.625
$ws.Cells.Item(122, 9).Value = 1344.2
$ws.Cells.Item(122, 11).Value = 4032.6
$ws.Cells.Item(122, 13).Value = -1582.6
$ws.Cells.Item(132, 8).Value = 2516.8096
$ws.Cells.Item(132, 9).Value = 2391.8462
$ws.Cells.Item(132, 10).Value = 2719.875
$ws.Cells.Item(132, 11).Value = 7175.5386
$ws.Cells.Item(132, 12).Value = 8159.625
$ws.Cells.Item(132, 13).Value = -4645.5386
$ws.Cells.Item(132, 14).Value = -13219.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1732.6666
$ws.Cells.Item(7, 9).Value = 1478.2
$ws.Cells.Item(7, 11).Value = 1478.2
$ws.Cells.Item(7, 13).Value = -1366.2
$ws.Cells.Item(40, 8).Value = 2891.5293
$ws.Cells.Item(40, 9).Value = 2650.9092
$ws.Cells.Item(40, 10).Value = 3332.6667
$ws.Cells.Item(40, 11).Value = 2650.9092
$ws.Cells.Item(40, 12).Value = 3332.6667
$ws.Cells.Item(40, 13).Value = -2514.9092
$ws.Cells.Item(40, 14).Value = -3604.6667
$ws.Cells.Item(61, 8).Value = 1747.1666
$ws.Cells.Item(61, 9).Value = 1523.2858
$ws.Cells.Item(61, 10).Value = 2060.6
$ws.Cells.Item(61, 11).Value = 1523.2858
$ws.Cells.Item(61, 12).Value = 2060.6
$ws.Cells.Item(61, 13).Value = -1321.2858
$ws.Cells.Item(61, 14).Value = -2464.6
$ws.Cells.Item(113, 8).Value = 1747.1666
$ws.Cells.Item(113, 9).Value = 1523.2858
$ws.Cells.Item(113, 10).Value = 2060.6
$ws.Cells.Item(113, 11).Value = 1523.2858
$ws.Cells.Item(113, 12).Value = 2060.6
$ws.Cells.Item(113, 13).Value = 646.7141999999999
$ws.Cells.Item(113, 14).Value = -6400.6
$ws.Cells.Item(122, 8).Value = 22730564
$ws.Cells.Item(122, 9).Value = 35717316
$ws.Cells.Item(122, 10).Value = 3749.5
$ws.Cells.Item(122, 11).Value = 107151948
$ws.Cells.Item(122, 12).Value = 11248.5
$ws.Cells.Item(122, 13).Value = -107149498
$ws.Cells.Item(122, 14).Value = -16148.5
$ws.Cells.Item(126, 8).Value = 1732.6666
$ws.Cells.Item(126, 9).Value = 1478.2
$ws.Cells.Item(126, 11).Value = 4434.6
$ws.Cells.Item(126, 13).Value = -1964.6
$ws.Cells.Item(132, 8).Value = 147257.28
$ws.Cells.Item(132, 9).Value = 5649.5
$ws.Cells.Item(132, 11).Value = 16948.5
$ws.Cells.Item(132, 13).Value = -14418.5
$ws.Cells.Item(136, 8).Value = 51004
$ws.Cells.Item(136, 9).Value = 51004
$ws.Cells.Item(136, 11).Value = 153012
$ws.Cells.Item(136, 13).Value = -150462
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 13005202
$ws.Cells.Item(122, 9).Value = 14449891
$ws.Cells.Item(122, 11).Value = 43349673
$ws.Cells.Item(122, 13).Value = -43347223
$ws.Cells.Item(126, 8).Value = 62501460
$ws.Cells.Item(126, 9).Value = 66668124
$ws.Cells.Item(126, 10).Value = 1480
$ws.Cells.Item(126, 11).Value = 200004372
$ws.Cells.Item(126, 12).Value = 4440
$ws.Cells.Item(126, 13).Value = -200001902
$ws.Cells.Item(126, 14).Value = -9380
$ws.Cells.Item(132, 8).Value = 5892.385
$ws.Cells.Item(132, 9).Value = 6100.5
$ws.Cells.Item(132, 11).Value = 18301.5
$ws.Cells.Item(132, 13).Value = -15771.5
$ws.Cells.Item(136, 8).Value = 1199.3334
$ws.Cells.Item(136, 9).Value = 473.25
$ws.Cells.Item(136, 10).Value = 1780.2
$ws.Cells.Item(136, 11).Value = 1419.75
$ws.Cells.Item(136, 12).Value = 5340.6
$ws.Cells.Item(136, 13).Value = 1130.25
$ws.Cells.Item(136, 14).Value = -10440.6
